$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "Lucene"

# Update tab ratio (cosmetic; may not persist)
$wb.Windows.Item(1).TabRatio = 986

# Reorder header row: C1=F1, D1=Accuracy, E1=Precision, F1=Recall
$ws.Range("C1").Value = "F1"
$ws.Range("D1").Value = "Accuracy"
$ws.Range("E1").Value = "Precision"
$ws.Range("F1").Value = "Recall "

# Update data rows 2-11 (models x configurations, with F1 moved before Accuracy)
$ws.Range("A2").Value = "Logistic Regression"
$ws.Range("B2").Value = "Count Vectorizer"
$ws.Range("C2").Value = 88.32
$ws.Range("D2").Value = 82.46
$ws.Range("E2").Value = 84.42
$ws.Range("F2").Value = 92.86

$ws.Range("A3").Value = "Multinomial Naive Bayes"
$ws.Range("B3").Value = "Count Vectorizer"
$ws.Range("C3").Value = 88.01
$ws.Range("D3").Value = 81.11
$ws.Range("E3").Value = 80.4
$ws.Range("F3").Value = 97.45

$ws.Range("A4").Value = "Support Vector Machines"
$ws.Range("B4").Value = "Count Vectorizer"
$ws.Range("C4").Value = 83.86
$ws.Range("D4").Value = 72.7
$ws.Range("E4").Value = 72.76
$ws.Range("F4").Value = 99.79

$ws.Range("A5").Value = "Decision Tree"
$ws.Range("B5").Value = "Count Vectorizer"
$ws.Range("C5").Value = 87.8
$ws.Range("D5").Value = 81.35
$ws.Range("E5").Value = 82.76
$ws.Range("F5").Value = 93.76

$ws.Range("A6").Value = "Random Forest"
$ws.Range("B6").Value = "Count Vectorizer"
$ws.Range("C6").Value = 88.33
$ws.Range("D6").Value = 82.09
$ws.Range("E6").Value = 83.02
$ws.Range("F6").Value = 94.55

$ws.Range("A7").Value = "Logistic Regression"
$ws.Range("B7").Value = "CV + tfidf"
$ws.Range("C7").Value = 88.91
$ws.Range("D7").Value = 82.66
$ws.Range("E7").Value = 81.66
$ws.Range("F7").Value = 97.8

$ws.Range("A8").Value = "Multinomial Naive Bayes"
$ws.Range("B8").Value = "CV + tfidf"
$ws.Range("C8").Value = 88.44
$ws.Range("D8").Value = 81.48
$ws.Range("E8").Value = 79.75
$ws.Range("F8").Value = 99.54

$ws.Range("A9").Value = "Support Vector Machines"
$ws.Range("B9").Value = "CV + tfidf"
$ws.Range("C9").Value = 83.94
$ws.Range("D9").Value = 72.78
$ws.Range("E9").Value = 72.78
$ws.Range("F9").Value = 100

$ws.Range("A10").Value = "Decision Tree"
$ws.Range("B10").Value = "CV + tfidf"
$ws.Range("C10").Value = 87.8
$ws.Range("D10").Value = 81.43
$ws.Range("E10").Value = 82.64
$ws.Range("F10").Value = 93.96

$ws.Range("A11").Value = "Random Forest"
$ws.Range("B11").Value = "CV + tfidf"
$ws.Range("C11").Value = 88.23
$ws.Range("D11").Value = 81.96
$ws.Range("E11").Value = 83.04
$ws.Range("F11").Value = 94.37

# Update selection to match target (F11)
$ws.Range("F11").Select()
